$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 11
$ws.Range("H11").Value = 1921
$ws.Range("I11").Value = 1921
$ws.Range("K11").Value = 1921
$ws.Range("M11").Value = -1781

# Row 12
$ws.Range("H12").Value = 175
$ws.Range("I12").Value = 133.33333
$ws.Range("J12").Value = 300
$ws.Range("K12").Value = 133.33333
$ws.Range("L12").Value = 300
$ws.Range("M12").Value = 36.66667000000001
$ws.Range("N12").Value = -640

# Row 17
$ws.Range("H17").Value = 1852164.6
$ws.Range("J17").Value = 1852164.6
$ws.Range("L17").Value = 5556493.800000001
$ws.Range("N17").Value = -5556829.800000001

# Row 43
$ws.Range("H43").Value = 938.5833
$ws.Range("I43").Value = 759.6667
$ws.Range("J43").Value = 998.2222
$ws.Range("K43").Value = 759.6667
$ws.Range("L43").Value = 998.2222
$ws.Range("M43").Value = -690.6667
$ws.Range("N43").Value = -1136.2222

# Row 132
$ws.Range("H132").Value = 2767.0952
$ws.Range("I132").Value = 2767.0952
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 8301.285600000001
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -5771.285600000001
$ws.Range("N132").ClearContents()

# Row 138
$ws.Range("H138").Value = 11630652
$ws.Range("I138").Value = 1310.1957
$ws.Range("J138").Value = 25004396
$ws.Range("K138").Value = 3930.5871
$ws.Range("L138").Value = 75013188
$ws.Range("M138").Value = 1209.4129
$ws.Range("N138").Value = -75023468

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3502.7341
$ws.Range("I32").Value = 2158.3242
$ws.Range("J32").Value = 23400
$ws.Range("K32").Value = 2158.3242
$ws.Range("L32").Value = 23400
$ws.Range("M32").Value = -1871.3242
$ws.Range("N32").Value = -23974

# Row 61
$ws.Range("H61").Value = 5102.06
$ws.Range("I61").Value = 5495.0967
$ws.Range("J61").Value = 4460.7896
$ws.Range("K61").Value = 5495.0967
$ws.Range("L61").Value = 4460.7896
$ws.Range("M61").Value = -5283.0967
$ws.Range("N61").Value = -4884.7896

# Row 74
$ws.Range("H74").Value = 1327.6666
$ws.Range("I74").Value = 1437.2667
$ws.Range("J74").Value = 1053.6666
$ws.Range("K74").Value = 1437.2667
$ws.Range("L74").Value = 1053.6666
$ws.Range("M74").Value = -563.2666999999999
$ws.Range("N74").Value = -2801.6666

# Row 77
$ws.Range("H77").Value = 1327.6666
$ws.Range("I77").Value = 1437.2667
$ws.Range("J77").Value = 1053.6666
$ws.Range("K77").Value = 7186.3335
$ws.Range("L77").Value = 5268.333000000001
$ws.Range("M77").Value = -2818.3335
$ws.Range("N77").Value = -14004.333

# Row 102
$ws.Range("H102").Value = 2521.2222
$ws.Range("I102").Value = 2397.5
$ws.Range("K102").Value = 2397.5
$ws.Range("M102").Value = -775.5

# Row 110
$ws.Range("H110").Value = 76007.375
$ws.Range("I110").Value = 120466.6
$ws.Range("J110").Value = 1908.6666
$ws.Range("K110").Value = 120466.6
$ws.Range("L110").Value = 1908.6666
$ws.Range("M110").Value = -118421.6
$ws.Range("N110").Value = -5998.6666

# Row 136
$ws.Range("H136").Value = 5102.06
$ws.Range("I136").Value = 5495.0967
$ws.Range("J136").Value = 4460.7896
$ws.Range("K136").Value = 16485.2901
$ws.Range("L136").Value = 13382.3688
$ws.Range("M136").Value = -13935.2901
$ws.Range("N136").Value = -18482.3688

# Row 141
$ws.Range("H141").Value = 74881.336
$ws.Range("J141").Value = 74881.336
$ws.Range("L141").Value = 74881.336
$ws.Range("N141").Value = -85241.336

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 1995.4
$ws.Range("I86").Value = 2227.4546
$ws.Range("J86").Value = 1711.7778
$ws.Range("K86").Value = 2227.4546
$ws.Range("L86").Value = 1711.7778
$ws.Range("M86").Value = -1104.4546
$ws.Range("N86").Value = -3957.7778

# Row 89
$ws.Range("H89").Value = 1995.4
$ws.Range("I89").Value = 2227.4546
$ws.Range("J89").Value = 1711.7778
$ws.Range("K89").Value = 11137.273
$ws.Range("L89").Value = 8558.889000000001
$ws.Range("M89").Value = -5521.273000000001
$ws.Range("N89").Value = -19790.889

# Row 105
$ws.Range("H105").Value = 1544.4375
$ws.Range("I105").Value = 1530
$ws.Range("K105").Value = 1530
$ws.Range("M105").Value = 217

# Row 122
$ws.Range("H122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()

# Row 135
$ws.Range("H135").Value = 44980
$ws.Range("J135").Value = 44980
$ws.Range("L135").Value = 44980
$ws.Range("N135").Value = -55120

# Row 137
$ws.Range("H137").Value = 30000
$ws.Range("J137").Value = 30000
$ws.Range("L137").Value = 30000
$ws.Range("N137").Value = -40200

# Row 140
$ws.Range("H140").Value = 61026
$ws.Range("J140").Value = 61026
$ws.Range("L140").Value = 61026
$ws.Range("N140").Value = -71386

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 6
$ws.Range("H6").Value = 10352111
$ws.Range("I6").Value = 10352111
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 10352111
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -10351998
$ws.Range("N6").ClearContents()

# Row 58
$ws.Range("H58").Value = 1279.8422
$ws.Range("I58").Value = 1322.8334
$ws.Range("J58").Value = 1206.1428
$ws.Range("K58").Value = 1322.8334
$ws.Range("L58").Value = 1206.1428
$ws.Range("M58").Value = -1119.8334
$ws.Range("N58").Value = -1612.1428

# Row 134
$ws.Range("H134").Value = 5353.129
$ws.Range("I134").Value = 5967.231
$ws.Range("K134").Value = 17901.693
$ws.Range("M134").Value = -15366.693

# Row 136
$ws.Range("H136").Value = 1279.8422
$ws.Range("I136").Value = 1322.8334
$ws.Range("J136").Value = 1206.1428
$ws.Range("K136").Value = 3968.5002
$ws.Range("L136").Value = 3618.4284
$ws.Range("M136").Value = -1418.5002
$ws.Range("N136").Value = -8718.428400000001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 59
$ws.Range("H59").Value = 3626.6667
$ws.Range("I59").Value = 480
$ws.Range("J59").Value = 4256
$ws.Range("K59").Value = 1440
$ws.Range("L59").Value = 12768
$ws.Range("M59").Value = -900
$ws.Range("N59").Value = -13848

# Row 81
$ws.Range("H81").Value = 4359.2856
$ws.Range("J81").Value = 4359.2856
$ws.Range("L81").Value = 13077.8568
$ws.Range("N81").Value = -15323.8568

# Row 84
$ws.Range("H84").Value = 4359.2856
$ws.Range("J84").Value = 4359.2856
$ws.Range("L84").Value = 39233.5704
$ws.Range("N84").Value = -50465.5704

# Row 112
$ws.Range("H112").Value = 1644.2858
$ws.Range("J112").Value = 1680.7693
$ws.Range("L112").Value = 5042.3079
$ws.Range("N112").Value = -7258.3079

# Row 131
$ws.Range("H131").Value = 4645.3667
$ws.Range("I131").Value = 506.66666
$ws.Range("J131").Value = 5680.0415
$ws.Range("K131").Value = 1519.99998
$ws.Range("L131").Value = 17040.1245
$ws.Range("M131").Value = 3520.00002
$ws.Range("N131").Value = -27120.1245

# Row 136
$ws.Range("H136").Value = 50843.375
$ws.Range("I136").Value = 2408.3572
$ws.Range("J136").Value = 389888.5
$ws.Range("K136").Value = 7225.071599999999
$ws.Range("L136").Value = 1169665.5
$ws.Range("M136").Value = -2125.071599999999
$ws.Range("N136").Value = -1179865.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 3899.6667
$ws.Range("I80").Value = 3850
$ws.Range("K80").Value = 3850
$ws.Range("M80").Value = -2852

# Row 83
$ws.Range("H83").Value = 3899.6667
$ws.Range("I83").Value = 3850
$ws.Range("K83").Value = 19250
$ws.Range("M83").Value = -14258

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1164.4286
$ws.Range("I22").Value = 790
$ws.Range("J22").Value = 2100.5
$ws.Range("K22").Value = 790
$ws.Range("L22").Value = 2100.5
$ws.Range("M22").Value = -495
$ws.Range("N22").Value = -2690.5

# Row 27
$ws.Range("H27").Value = 1164.4286
$ws.Range("I27").Value = 790
$ws.Range("J27").Value = 2100.5
$ws.Range("K27").Value = 790
$ws.Range("L27").Value = 2100.5
$ws.Range("M27").Value = -683
$ws.Range("N27").Value = -2314.5

# Row 123
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 86235.836
$ws.Range("I122").Value = 169555
$ws.Range("J122").Value = 2916.6667
$ws.Range("K122").Value = 508665
$ws.Range("L122").Value = 8750.000100000001
$ws.Range("M122").Value = -506215
$ws.Range("N122").Value = -13650.0001

# Row 136
$ws.Range("H136").Value = 7250750.5
$ws.Range("I136").Value = 15152218
$ws.Range("K136").Value = 45456654
$ws.Range("M136").Value = -45454104
